# Auto-generated Excel COM-interop script applying numeric updates
# to the Mateus_Profits leve-crafting profit workbook (columns H-N per row),
# matching the target OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2180
$ws.Range("I2").Value = 950
$ws.Range("K2").Value = 950
$ws.Range("M2").Value = -837

$ws.Range("H17").Value = 6252355.5
$ws.Range("J17").Value = 6669112.5
$ws.Range("L17").Value = 20007337.5
$ws.Range("N17").Value = -20007673.5

$ws.Range("H28").Value = 620.625
$ws.Range("I28").Value = 761.5
$ws.Range("K28").Value = 761.5
$ws.Range("M28").Value = -276.5

$ws.Range("H43").Value = 8100
$ws.Range("I43").Value = 8500
$ws.Range("K43").Value = 8500
$ws.Range("M43").Value = -8431

$ws.Range("H132").Value = 8636.925999999999
$ws.Range("I132").Value = 1327.88
$ws.Range("K132").Value = 3983.64
$ws.Range("M132").Value = -1453.64

$ws.Range("H137").Value = 3630.087
$ws.Range("I137").Value = 2507.375
$ws.Range("J137").Value = 4228.8667
$ws.Range("K137").Value = 7522.125
$ws.Range("L137").Value = 12686.6001
$ws.Range("M137").Value = -4972.125
$ws.Range("N137").Value = -17786.6001

$ws.Range("H138").Value = 5582.8184
$ws.Range("J138").Value = 5983.3687
$ws.Range("L138").Value = 17950.1061
$ws.Range("N138").Value = -28230.1061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 230
$ws.Range("I4").Value = 183.22728
$ws.Range("K4").Value = 183.22728
$ws.Range("M4").Value = -67.22728000000001

$ws.Range("H32").Value = 9663.5
$ws.Range("I32").Value = 8670.743
$ws.Range("J32").Value = 17407
$ws.Range("K32").Value = 8670.743
$ws.Range("L32").Value = 17407
$ws.Range("M32").Value = -8383.743
$ws.Range("N32").Value = -17981

$ws.Range("H45").Value = 147284.28
$ws.Range("I45").Value = 225053.33
$ws.Range("K45").Value = 225053.33
$ws.Range("M45").Value = -224676.33

$ws.Range("H63").Value = 5845.5
$ws.Range("I63").Value = 3166.3333
$ws.Range("K63").Value = 3166.3333
$ws.Range("M63").Value = -2480.3333

$ws.Range("H66").Value = 5845.5
$ws.Range("I66").Value = 3166.3333
$ws.Range("K66").Value = 15831.6665
$ws.Range("M66").Value = -12399.6665

$ws.Range("H122").Value = 2354.4736
$ws.Range("I122").Value = 2196.2354
$ws.Range("J122").Value = 3699.5
$ws.Range("K122").Value = 6588.706200000001
$ws.Range("L122").Value = 11098.5
$ws.Range("M122").Value = -4138.706200000001
$ws.Range("N122").Value = -15998.5

$ws.Range("H132").Value = 3390.238
$ws.Range("I132").Value = 3390.238
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10170.714
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7640.714
$ws.Range("N132").Value = ""

$ws.Range("H134").Value = 163748.5
$ws.Range("J134").Value = 163748.5
$ws.Range("L134").Value = 163748.5
$ws.Range("N134").Value = -173888.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7313.9443
$ws.Range("I31").Value = 3499.5
$ws.Range("K31").Value = 3499.5
$ws.Range("M31").Value = -3204.5

$ws.Range("H34").Value = 7313.9443
$ws.Range("I34").Value = 3499.5
$ws.Range("K34").Value = 3499.5
$ws.Range("M34").Value = -3297.5

$ws.Range("H58").Value = 4702.6113
$ws.Range("I58").Value = 2443.3044
$ws.Range("K58").Value = 2443.3044
$ws.Range("M58").Value = -2240.3044

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""

$ws.Range("H133").Value = 46456.125
$ws.Range("I133").Value = 30000
$ws.Range("K133").Value = 30000
$ws.Range("M133").Value = -27470

$ws.Range("H134").Value = 4629.5483
$ws.Range("I134").Value = 3633.5186
$ws.Range("J134").Value = 11352.75
$ws.Range("K134").Value = 10900.5558
$ws.Range("L134").Value = 34058.25
$ws.Range("M134").Value = -8365.5558
$ws.Range("N134").Value = -39128.25

$ws.Range("H136").Value = 4702.6113
$ws.Range("I136").Value = 2443.3044
$ws.Range("K136").Value = 7329.9132
$ws.Range("M136").Value = -4779.9132

$ws.Range("H140").Value = 123256.2
$ws.Range("J140").Value = 123256.2
$ws.Range("L140").Value = 123256.2
$ws.Range("N140").Value = -133616.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 667059.3
$ws.Range("I7").Value = 667059.3
$ws.Range("K7").Value = 2001177.9
$ws.Range("M7").Value = -2001065.9

$ws.Range("H15").Value = 93880.625
$ws.Range("I15").Value = 165
$ws.Range("K15").Value = 495
$ws.Range("M15").Value = -355

$ws.Range("H16").Value = 1550
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 7500
$ws.Range("M16").Value = -7327

$ws.Range("H20").Value = 2000
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 6000
$ws.Range("N20").Value = -6454

$ws.Range("H21").Value = 125200
$ws.Range("J21").Value = 250000
$ws.Range("L21").Value = 750000
$ws.Range("N21").Value = -750346

$ws.Range("H22").Value = 1336.45
$ws.Range("J22").Value = 1376.9445
$ws.Range("L22").Value = 4130.833500000001
$ws.Range("N22").Value = -4468.833500000001

$ws.Range("H26").Value = 224766.78
$ws.Range("J26").Value = 3799.3333
$ws.Range("L26").Value = 11397.9999
$ws.Range("N26").Value = -11973.9999

$ws.Range("H27").Value = 1336.45
$ws.Range("J27").Value = 1376.9445
$ws.Range("L27").Value = 4130.833500000001
$ws.Range("N27").Value = -4334.833500000001

$ws.Range("H32").Value = 525
$ws.Range("J32").Value = 525
$ws.Range("L32").Value = 1575
$ws.Range("N32").Value = -2141

$ws.Range("H33").Value = 383.6154
$ws.Range("J33").Value = 805.5
$ws.Range("L33").Value = 4833
$ws.Range("N33").Value = -5399

$ws.Range("H34").Value = 100387.5
$ws.Range("I34").Value = 395.83334
$ws.Range("J34").Value = 250375
$ws.Range("K34").Value = 1187.50002
$ws.Range("L34").Value = 751125
$ws.Range("M34").Value = -1103.50002
$ws.Range("N34").Value = -751293

$ws.Range("H63").Value = 250002220
$ws.Range("J63").Value = 2888
$ws.Range("L63").Value = 8664
$ws.Range("N63").Value = -10162

$ws.Range("H66").Value = 250002220
$ws.Range("J66").Value = 2888
$ws.Range("L66").Value = 25992
$ws.Range("N66").Value = -33480

$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 50000
$ws.Range("K74").Value = 150000
$ws.Range("M74").Value = -148939

$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 50000
$ws.Range("K77").Value = 450000
$ws.Range("M77").Value = -444696

$ws.Range("H92").Value = 348
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws.Range("H139").Value = 5238.2
$ws.Range("I139").Value = 1970.5862
$ws.Range("J139").Value = 99999
$ws.Range("K139").Value = 5911.7586
$ws.Range("L139").Value = 299997
$ws.Range("M139").Value = -771.7586000000001
$ws.Range("N139").Value = -310277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 254958.33
$ws.Range("I21").Value = 303950
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 303950
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -303777
$ws.Range("N21").Value = -10346

$ws.Range("H30").Value = 254958.33
$ws.Range("I30").Value = 303950
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 303950
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = -303845
$ws.Range("N30").Value = -10210

$ws.Range("H80").Value = 2553.8948
$ws.Range("J80").Value = 2389.625
$ws.Range("L80").Value = 2389.625
$ws.Range("N80").Value = -4385.625

$ws.Range("H83").Value = 2553.8948
$ws.Range("J83").Value = 2389.625
$ws.Range("L83").Value = 11948.125
$ws.Range("N83").Value = -21932.125

$ws.Range("H102").Value = 3712.4707
$ws.Range("I102").Value = 2740.8667
$ws.Range("K102").Value = 2740.8667
$ws.Range("M102").Value = -1118.8667

$ws.Range("H135").Value = 56308.34
$ws.Range("J135").Value = 56308.34
$ws.Range("L135").Value = 56308.34
$ws.Range("N135").Value = -66448.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4035.5334
$ws.Range("I40").Value = 3704.25
$ws.Range("J40").Value = 4414.143
$ws.Range("K40").Value = 3704.25
$ws.Range("L40").Value = 4414.143
$ws.Range("M40").Value = -3568.25
$ws.Range("N40").Value = -4686.143

$ws.Range("H46").Value = 2900
$ws.Range("I46").Value = 3049.6667
$ws.Range("J46").Value = 2451
$ws.Range("K46").Value = 3049.6667
$ws.Range("L46").Value = 2451
$ws.Range("M46").Value = -2861.6667
$ws.Range("N46").Value = -2827

$ws.Range("H130").Value = 38964
$ws.Range("J130").Value = 38964
$ws.Range("L130").Value = 38964
$ws.Range("N130").Value = -49004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 481.32
$ws.Range("J113").Value = 988
$ws.Range("L113").Value = 2964
$ws.Range("N113").Value = -7304

$ws.Range("H133").Value = 142499.5
$ws.Range("I133").Value = 40000
$ws.Range("K133").Value = 40000
$ws.Range("M133").Value = -34940
